# "Cap nhat Phieu van dap.xlsx" - update completion levels ("Muc do hoan thanh")
# for several checklist rows on the PhieuVanDap sheet. The dependent formula
# cells (F/G/H columns) recalculate automatically from these inputs.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PhieuVanDap")

# 1.1 "Dang ky tai khoan thanh vien thuong (Member)" -> hoan thanh 100%
$ws.Range("D13").Value = 1

# 3.1 "Xem va duyet tin rao vat" -> hoan thanh 100%
$ws.Range("D43").Value = 1

# 3.2 "Them, Xoa tin rao vat" -> hoan thanh 100%
$ws.Range("D44").Value = 1

# 3.3 (muc con cua phan 3) -> hoan thanh 50%
$ws.Range("D45").Value = 0.5

# 4.1 "Quan ly tin rao vat" -> hoan thanh 100%
$ws.Range("D47").Value = 1
